# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# values on the active worksheet with the latest scraped figures.
#
# Price values are stored as plain text in the workbook (the source feed
# renders them pre-formatted, e.g. "25.965.17" or "1.000", which Excel
# would otherwise auto-convert to a number). Briefly switching the cell to
# the Text number format while assigning the value - then restoring the
# plain "Normal" cell style - keeps the literal text intact without
# leaving any visible formatting change behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '25.965.17'; E = '  +0.74%  ' },
    @{ Row = 3; D = '1.748.22'; E = '  -0.20%  ' },
    @{ Row = 4; D = '1.000'; E = '  +0.12%  ' },
    @{ Row = 5; D = '236.03'; E = '  -0.37%  ' },
    @{ Row = 6; D = '0.9995'; E = '  +0.06%  ' },
    @{ Row = 7; D = '0.5212'; E = '  +2.95%  ' },
    @{ Row = 8; D = '0.2848'; E = '  +4.91%  ' },
    @{ Row = 9; D = '39.32'; E = '  -3.68%  ' },
    @{ Row = 10; D = '0.06141'; E = '  -1.15%  ' },
    @{ Row = 11; D = '1.755.93'; E = '  +0.17%  ' },
    @{ Row = 12; D = '0.07033'; E = '  +1.40%  ' },
    @{ Row = 13; D = '15.50'; E = '  -1.01%  ' },
    @{ Row = 14; D = '0.6464'; E = '  +3.49%  ' },
    @{ Row = 15; D = '4.529'; E = '  +0.83%  ' },
    @{ Row = 16; D = '77.51'; E = '  -1.27%  ' },
    @{ Row = 17; D = '0.9984'; E = '  -0.05%  ' },
    @{ Row = 18; D = '0.9985'; E = '  -0.09%  ' },
    @{ Row = 19; D = '25.973.29'; E = '  +0.68%  ' },
    @{ Row = 20; D = '11.51'; E = '  -1.66%  ' },
    @{ Row = 21; D = '0.000006619'; E = '  -1.80%  ' },
    @{ Row = 22; D = '1.976.60'; E = '  +0.15%  ' },
    @{ Row = 23; D = '4.167'; E = '  +2.73%  ' },
    @{ Row = 24; D = '8.650'; E = '  +4.68%  ' },
    @{ Row = 25; D = '5.159'; E = '  -0.66%  ' },
    @{ Row = 26; D = '139.04'; E = '  +1.47%  ' },
    @{ Row = 27; D = '1.497'; E = '  +3.05%  ' },
    @{ Row = 28; D = '1.845'; E = '  +2.62%  ' },
    @{ Row = 29; D = $null; E = '  -0.42%  ' },
    @{ Row = 30; D = '103.19'; E = '  +0.19%  ' },
    @{ Row = 31; D = '0.08315'; E = '  +0.39%  ' },
    @{ Row = 32; D = '3.656'; E = '  -2.40%  ' },
    @{ Row = 33; D = '3.435'; E = '  +0.05%  ' },
    @{ Row = 34; D = '0.04471'; E = '  +1.87%  ' },
    @{ Row = 35; D = '2.613'; E = '  -1.25%  ' },
    @{ Row = 36; D = '0.9864'; E = '  -2.27%  ' },
    @{ Row = 37; D = '0.6105'; E = '  +0.92%  ' },
    @{ Row = 38; D = '2.691'; E = '  +0.25%  ' },
    @{ Row = 39; D = '0.01594'; E = '  +2.31%  ' },
    @{ Row = 40; D = '1.955'; E = '  -0.64%  ' },
    @{ Row = 41; D = '0.9983'; E = '  -0.07%  ' },
    @{ Row = 42; D = '101.03'; E = '  -1.00%  ' },
    @{ Row = 43; D = '0.3876'; E = '  +0.90%  ' },
    @{ Row = 44; D = '0.7368'; E = '  -2.25%  ' },
    @{ Row = 45; D = '5.064'; E = '  +4.23%  ' },
    @{ Row = 46; D = '0.05476'; E = '  -0.34%  ' },
    @{ Row = 47; D = '6.347'; E = '  +6.16%  ' },
    @{ Row = 48; D = '0.1120'; E = '  +3.07%  ' },
    @{ Row = 49; D = '53.01'; E = '  +1.00%  ' },
    @{ Row = 50; D = '30.09'; E = '  -0.71%  ' },
    @{ Row = 51; D = '7.630'; E = '  +1.43%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $priceCell = $ws.Cells.Item($u.Row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.D
        $priceCell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

Write-Output "Updated $($updates.Count) rows"
